$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $true, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# --- 1) Remove the trailing manual line break after "N0 ... lymph nodes" bullet ---
# (Leave the preceding runs untouched: only delete the single <w:br/> character.)
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t -like "N0 cancers have not spread to the lymph nodes*") {
        $e = $p.Range.End
        $brk = $d.Range($e - 2, $e - 1)
        if ($brk.Text -eq [string][char]11) {
            $brk.Delete() | Out-Null
        }
        break
    }
}

# --- 2) Laparoscopy section ---
# Remove the "A laparoscopy is performed under a general anesthetic." paragraph entirely
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "A laparoscopy is performed under a general anesthetic.*") {
        $p.Range.Delete() | Out-Null
        break
    }
}

# Add a new bullet "General anesthetic" right before "Several incisions 1/4” long"
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Several incisions 1/4*long*") {
        $p.Range.InsertParagraphBefore() | Out-Null
        $newP = $d.Paragraphs.Item($i)
        $newP.Range.Text = "General anesthetic"
        break
    }
}

# Telescope sentence
Replace-Text "A telescope is inserted to look inside the abdominal cavity." "A telescope is used to examine the abdomen"

# --- Simple wording tweaks ---
Replace-Text "Placed at the beginning of each dose" "Placed for each dose"
Replace-Text "Removed that day at the end of treatment" "Removed that day"
Replace-Text "May shower within 24 hrs" "May shower in 24 hrs"
Replace-Text "Placed underneath the skin below the right collarbone" "Implanted under the skin"
Replace-Text "Incision in the neck (1/4”)" "Neck incision (1/4”)"

# --- 9) Remove the empty "First Paragraph" (just a manual line break) before
#        "Critical to good communication..." ---
$brkOnly = [string][char]11 + [string][char]13
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Style.NameLocal -eq "First Paragraph" -and $p.Range.Text -eq $brkOnly) {
        $p.Range.Delete() | Out-Null
        break
    }
}

# --- Remaining wording tweaks ---
Replace-Text "Critical to good communication with your cancer care team" "Critical to good communication with your care team"
Replace-Text "Important to reduce the risk of complications from cancer treatment" "Reduces risk of complications from treatment"
Replace-Text "Working hard enough that you can’t carry a conversation" "Working hard enough that you can’t converse"
Replace-Text "Start slow an build up" "Start slowly and build up"
Replace-Text "Smoking makes it more difficult to get through cancer treatment" "Smoking makes cancer treatment more difficult"
Replace-Text "American Lung Asssociation fredomfromsmoking.org" "American Lung Assn fredomfromsmoking.org"
Replace-Text "1:1 Smoking Cessation Counseling Clinics (Metro Charlotte)" "1:1 Smoking Cessation Counseling (Metro Charlotte)"
